$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Mirror column G (rows 7-29) into a new column H, preserving both value and
# direct formatting (wrapText / font style) for every cell in that range,
# including the currently-empty ones.
for ($r = 7; $r -le 29; $r++) {
    $src = $ws.Range("G$r")
    $dst = $ws.Range("H$r")
    $src.Copy()
    $dst.PasteSpecial(-4122)   # xlPasteFormats
    $dst.Value2 = $src.Value2
}
$excel.CutCopyMode = $false

# A7 previously carried a stray "applyFont" style left over from a copy/paste;
# strip it back to the default formatting.
$ws.Range("A7").ClearFormats()

# A8 was an empty cell that only existed to hold that same stray style -
# remove it entirely now that the style is gone.
$ws.Range("A8").Clear()

# Update the selection to reflect where editing continued.
[void]$ws.Range("H4").Select()
